$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new entry for "CS intoduction Lecture 13"
$ws.Range("A14").Value = 45813
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 50
$ws.Range("D14").Value = 17
$ws.Range("E14").Value = 14
$ws.Range("F14").Value = "CS intoduction Lecture 13"

# Row 15: new entry for "CS intoduction Lecture 14"
$ws.Range("A15").Value = 45813
$ws.Range("B15").Value = 17
$ws.Range("C15").Value = 15
$ws.Range("F15").Value = "CS intoduction Lecture 14"

# Match the date number formatting used by the rest of column A (copy from A13)
# without introducing new style/numFmt entries.
$ws.Range("A13").Copy()
$ws.Range("A14:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F15").Select()
